# Append three new M2 data rows (225-227) to the Romania_M2 sheet,
# mirroring the formatting of the last existing data row (224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: datetime(serial), symbol, open, high, low, close, volume
$newRows = @(
    @(45108.41666666666, "ECONOMICS:ROM2", 622009300000, 622009300000, 622009300000, 622009300000, 0),
    @(45139.41666666666, "ECONOMICS:ROM2", 626029700000, 626029700000, 626029700000, 626029700000, 0),
    @(45170.41666666666, "ECONOMICS:ROM2", 640762400000, 640762400000, 640762400000, 640762400000, 0)
)

$lastRow = 224
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $targetRow = $startRow + $i
    $rowValues = $newRows[$i]

    # Copy the formatting (styles/number formats) of the last data row down
    # to the new row so the appended rows look identical to the existing ones.
    $ws.Range("A" + $lastRow + ":G" + $lastRow).Copy()
    $ws.Range("A" + $targetRow + ":G" + $targetRow).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($targetRow, 1).Value2 = $rowValues[0]
    $ws.Cells.Item($targetRow, 2).Value2 = $rowValues[1]
    $ws.Cells.Item($targetRow, 3).Value2 = $rowValues[2]
    $ws.Cells.Item($targetRow, 4).Value2 = $rowValues[3]
    $ws.Cells.Item($targetRow, 5).Value2 = $rowValues[4]
    $ws.Cells.Item($targetRow, 6).Value2 = $rowValues[5]
    $ws.Cells.Item($targetRow, 7).Value2 = $rowValues[6]
}

$excel.CutCopyMode = 0
